# Fixed up(enumsFinuras): return "DIA" in tablerows
#
# 1) B1 gets the header label "Unnamed: 1" (was an empty string cell).
# 2) Rows 3-8 are appended as blank data rows, each with a handful of
#    numeric values scattered in the JACQUARD block of columns.
# 3) Row 9 is appended as the month's totals row (A9 = 24, plus the same
#    JACQUARD numeric columns).
# All other cells in rows 3-9 are empty strings, matching the existing
# empty cells already present in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 75  # column BW

# --- B1 header -------------------------------------------------------
$ws.Cells.Item(1, 2).Value = "Unnamed: 1"

# --- Data for the new rows (1-based column numbers -> value) ---------
$rowData = @{
    3 = @{ 42 = 4; 44 = 4; 46 = 5; 48 = 4; 51 = 9 }
    4 = @{ 38 = 4; 48 = 4; 51 = 4 }
    5 = @{ 38 = 4; 44 = 4; 46 = 4; 51 = 8 }
    6 = @{ 38 = 4; 42 = 4; 46 = 4; 51 = 12 }
    7 = @{ 38 = 4; 42 = 4; 46 = 4; 51 = 12 }
    8 = @{ 38 = 4; 42 = 4; 46 = 4; 51 = 12 }
    9 = @{ 1 = 24; 38 = 4; 44 = 2; 46 = 6; 51 = 10 }
}

foreach ($r in 3..9) {
    $values = $rowData[$r]
    for ($c = 1; $c -le $lastCol; $c++) {
        if ($values.ContainsKey($c)) {
            $ws.Cells.Item($r, $c).Value = $values[$c]
        } else {
            $ws.Cells.Item($r, $c).Value = ""
        }
    }
}
